$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Iraq League")

function Swap-Rows {
    param($row1, $row2)

    for ($col = 2; $col -le 30; $col++) {
        $c1 = $ws.Cells.Item($row1, $col)
        $c2 = $ws.Cells.Item($row2, $col)
        $v1 = $c1.Value2
        $v2 = $c2.Value2
        $c1.Value2 = $v2
        $c2.Value2 = $v1
    }
}

Swap-Rows 17 18
Swap-Rows 58 59
Swap-Rows 78 79
Swap-Rows 135 136
Swap-Rows 219 220
